# Refresh PRAN yearly financial figures (Income Statement, Balance Sheet,
# and Cash Flow Statement) with updated source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRAN")

# Row 12: Research Development
$ws.Range("D12").Value = 4700
$ws.Range("E12").Value = 4000
$ws.Range("F12").Value = 6800
$ws.Range("G12").Value = 8700
$ws.Range("H12").Value = 10500
$ws.Range("I12").Value = 11400
$ws.Range("J12").Value = 6000

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 5800
$ws.Range("E17").Value = 5000
$ws.Range("F17").Value = 6200
$ws.Range("G17").Value = 12200
$ws.Range("H17").Value = 14700
$ws.Range("I17").Value = 8800
$ws.Range("J17").Value = 5500

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = -5700
$ws.Range("E18").Value = -4900
$ws.Range("F18").Value = -6100
$ws.Range("G18").Value = -12000
$ws.Range("H18").Value = -14400
$ws.Range("I18").Value = -8700
$ws.Range("J18").Value = -5400

# Row 20: Total Other Income/Expenses Net
$ws.Range("G20").Value = 7900
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 3200
$ws.Range("J20").Value = 1700

# Row 21: Earnings Before Interest And Taxes
$ws.Range("G21").Value = -4100
$ws.Range("H21").Value = -9400
$ws.Range("I21").Value = -5500
$ws.Range("J21").Value = -3700

# Row 23: Income Before Tax
$ws.Range("D23").Value = -5800
$ws.Range("E23").Value = -5300
$ws.Range("F23").Value = -5500
$ws.Range("G23").Value = -4200
$ws.Range("H23").Value = -9400
$ws.Range("I23").Value = -5500
$ws.Range("J23").Value = -3700

# Row 26: Income After Tax
$ws.Range("D26").Value = -5800
$ws.Range("E26").Value = -5300
$ws.Range("F26").Value = -5500
$ws.Range("G26").Value = -4200
$ws.Range("H26").Value = -9400
$ws.Range("I26").Value = -5500
$ws.Range("J26").Value = -3700

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = -5800
$ws.Range("E27").Value = -5300
$ws.Range("F27").Value = -5500
$ws.Range("G27").Value = -4200
$ws.Range("H27").Value = -9400
$ws.Range("I27").Value = -5500
$ws.Range("J27").Value = -3700

# Row 32: Other Items
$ws.Range("G32").Value = -7900
$ws.Range("H32").Value = -5000
$ws.Range("I32").Value = -3200
$ws.Range("J32").Value = -1700

# Row 33: Net Income
$ws.Range("D33").Value = -5800
$ws.Range("E33").Value = -5300
$ws.Range("F33").Value = -5500
$ws.Range("G33").Value = -4200
$ws.Range("H33").Value = -9400
$ws.Range("I33").Value = -5500
$ws.Range("J33").Value = -3700

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = -5800
$ws.Range("E35").Value = -5300
$ws.Range("F35").Value = -5500
$ws.Range("G35").Value = -4200
$ws.Range("H35").Value = -9400
$ws.Range("I35").Value = -5500
$ws.Range("J35").Value = -3700

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 10800
$ws.Range("E41").Value = 15500
$ws.Range("F41").Value = 20200
$ws.Range("G41").Value = 24700
$ws.Range("H41").Value = 24200
$ws.Range("I41").Value = 18900
$ws.Range("J41").Value = 4000

# Row 43: Net Receivables
$ws.Range("D43").Value = 2200
$ws.Range("E43").Value = 2100
$ws.Range("F43").Value = 3400
$ws.Range("G43").Value = 4600
$ws.Range("H43").Value = 5200
$ws.Range("I43").Value = 5000

# Row 46: Total Current Assets
$ws.Range("D46").Value = 13200
$ws.Range("E46").Value = 17900
$ws.Range("F46").Value = 23800
$ws.Range("G46").Value = 29500
$ws.Range("H46").Value = 29400
$ws.Range("I46").Value = 12000
$ws.Range("J46").Value = 5100

# Row 54: Total Assets
$ws.Range("D54").Value = 13200
$ws.Range("E54").Value = 17900
$ws.Range("F54").Value = 23900
$ws.Range("G54").Value = 29600
$ws.Range("H54").Value = 29500
$ws.Range("I54").Value = 12100
$ws.Range("J54").Value = 5200

# Row 57: Accounts Payable
$ws.Range("F57").Value = 1200

# Row 59: Other Current Liabilities
$ws.Range("I59").Value = 1700

# Row 60: Total Current Liabilities
$ws.Range("E60").Value = 1100
$ws.Range("G60").Value = 1900
$ws.Range("H60").Value = 2800

# Row 66: Total Liabilities
$ws.Range("E66").Value = 1100
$ws.Range("G66").Value = 1900
$ws.Range("H66").Value = 2800

# Row 72: Retained Earnings
$ws.Range("D72").Value = -90400
$ws.Range("E72").Value = -85100
$ws.Range("F72").Value = -81700
$ws.Range("G72").Value = -82900
$ws.Range("H72").Value = -78700
$ws.Range("I72").Value = -61800
$ws.Range("J72").Value = -57000

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 11400
$ws.Range("E76").Value = 16800
$ws.Range("F76").Value = 22200
$ws.Range("G76").Value = 27700
$ws.Range("H76").Value = 26700
$ws.Range("I76").Value = 9900
$ws.Range("J76").Value = 4000

# Row 81: Net Income
$ws.Range("D81").Value = -5800
$ws.Range("E81").Value = -5300
$ws.Range("F81").Value = -5500
$ws.Range("G81").Value = -4200
$ws.Range("H81").Value = -9400
$ws.Range("I81").Value = -5500
$ws.Range("J81").Value = -3700

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = -4400
$ws.Range("E89").Value = -4100
$ws.Range("F89").Value = -5100
$ws.Range("G89").Value = -7700
$ws.Range("H89").Value = -9600
$ws.Range("I89").Value = -5600
$ws.Range("J89").Value = -4800

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("G100").Value = 4800
$ws.Range("H100").Value = 24700
$ws.Range("I100").Value = 11000

# Row 101: Effect Of Exchange Rate Changes
$ws.Range("G101").Value = 3500

# Row 102: Change In Cash and Cash Equivalents
$ws.Range("D102").Value = -4700
$ws.Range("E102").Value = -4700
$ws.Range("F102").Value = -4500
$ws.Range("H102").Value = 14700
$ws.Range("I102").Value = 5500
